$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Drop the "anyone who cannot wear a mask..." sentence. The paragraph it
#    lives in stays (now empty) - only the run/text goes away.
# ---------------------------------------------------------------------------
$ok1 = $d.Content.Find.Execute(
    "anyone who cannot wear a mask should contact emma.sokell@ucd.ie",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 2)
if (-not $ok1) { throw "could not find the 'anyone who cannot...' sentence" }

# ---------------------------------------------------------------------------
# 2) That sentence's paragraph is immediately followed by a blank paragraph
#    that separated it from the "You should sign in/out..." paragraph. That
#    blank paragraph is removed (merged away) in the new revision.
# ---------------------------------------------------------------------------
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "") {
        $next = $p.Next()
        if ($next -ne $null -and $next.Range.Text -like "You should sign in/out*") {
            $target = $p
            break
        }
    }
}
if ($target -eq $null) { throw "could not locate the blank paragraph before 'You should sign in/out...'" }
$target.Range.Delete()

# ---------------------------------------------------------------------------
# 3) Drop the (hidden) "_GoBack" bookmark that sat inside the "You should
#    sign in/out..." paragraph - the surrounding text is untouched.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 4) Renumber item "13." (Bags, jackets...) to "12." - the numbered item
#    above it ("14 Computational experiments...") is about to disappear, so
#    this one shifts up to take its place.
# ---------------------------------------------------------------------------
$ok4 = $d.Content.Find.Execute(
    "13. Bags, jackets etc to be ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "12. Bags, jackets etc to be ", 2)
if (-not $ok4) { throw "could not find the '13. Bags, jackets...' item" }

# ---------------------------------------------------------------------------
# 5) Remove the whole "14 Computational experiments, introductory
#    exercised, data analysis..." paragraph.
# ---------------------------------------------------------------------------
$idx = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*Computational experiments*") {
        $idx = $i
        break
    }
}
if ($idx -eq 0) { throw "could not locate the 'Computational experiments' paragraph" }
$d.Paragraphs.Item($idx).Range.Delete()
